{"js": "const paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items[0].insertText(\"Answers: Rationalizing the denominator\", \"Replace\");\nparagraphs.items[1].insertText(\"Maximilian Volmar\", \"Replace\");\nparagraphs.items[3].insertText(\"Answers to questions relating to the guide on rationalizing the denominator.\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# Collapse the split runs in the Title, Author and Abstract paragraphs back\n# into a single run each, matching the original (un-split) text. We use\n# Find/Replace (ReplaceOne) scoped to the whole story but limited to a single\n# occurrence so later, unrelated text (\"... by Maximilian Volmar.\") in the\n# licensing section is left untouched.\n$d = $word.ActiveDocument\n\nfunction Replace-FirstOccurrence($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)\n}\n\nReplace-FirstOccurrence \"Answers: Rationalizing the denominator\" \"Answers: Rationalizing the denominator\"\nReplace-FirstOccurrence \"Maximilian Volmar\" \"Maximilian Volmar\"\nReplace-FirstOccurrence \"Answers to questions relating to the guide on rationalizing the denominator.\" \"Answers to questions relating to the guide on rationalizing the denominator.\"\n"}
